# AttendanceSheet001 - remove a student record ("Wain, ALEXANDRA") ahead of
# the next iteration of the course. Deleting the whole row shifts every
# subsequent roster row up by one, which matches the rest of the sheet
# shrinking by a single row (A1:O72 -> A1:O71) without touching any other
# student's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66 holds "Wain" / "ALEXANDRA" - drop the whole row.
$ws.Rows.Item(66).Delete()

# Leave the sheet's cursor where the instructor ended up after editing the
# bottom of the roster.
$ws.Range("H82").Select()
